$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Section: duplicatedFieldInDatatype2_negativeTest (rows 35-38) ----
$ws.Range("B35:E35").Merge()
$ws.Range("B35").Value = "Test doAnnualPremium duplicatedFieldInDatatype2_negativeTest"
$ws.Range("B35:E35").HorizontalAlignment = -4108

$ws.Range("B36").Value = "rate"
$ws.Range("C36").Value = "totalValue"
$ws.Range("D36").Value = "_res_.rate(1)"
$ws.Range("E36").Value = "_res_.rate"

$ws.Range("B37").Value = "rate"
$ws.Range("C37").Value = "totalValue"
$ws.Range("D37").Value = "_res_.rate"
$ws.Range("E37").Value = "_res_.rate"

$ws.Range("B38").Value = 0.05
$ws.Range("C38").Value = 5000
$ws.Range("C38").NumberFormat = '"$"#,##0.00'
$ws.Range("D38").Value = 250
$ws.Range("E38").Value = 240

# ---- Section: duplicatedFieldInSpr2_negativeTest (rows 43-46) ----
$ws.Range("B43:E43").Merge()
$ws.Range("B43").Value = "Test doAnnualPremiumSpr duplicatedFieldInSpr2_negativeTest"
$ws.Range("B43:E43").HorizontalAlignment = -4108

$ws.Range("B44").Value = "rate"
$ws.Range("C44").Value = "totalValue"
$ws.Range("D44").Value = "_res_.`$MonthYear"
$ws.Range("E44").Value = "_res_.`$MonthYear(-2)"

$ws.Range("B45").Value = "rate"
$ws.Range("C45").Value = "totalValue"
$ws.Range("D45").Value = "_res_.`$MonthYear"
$ws.Range("E45").Value = "_res_.`$MonthYear"

$ws.Range("B46").Value = 0.05
$ws.Range("C46").Value = 5000
$ws.Range("C46").NumberFormat = '"$"#,##0.00'
$ws.Range("D46").Value = 250
$ws.Range("E46").Value = 240

# ---- Section: duplicatedFieldInSpr3_negativeTest (rows 51-54) ----
$ws.Range("B51:E51").Merge()
$ws.Range("B51").Value = "Test doAnnualPremiumSpr duplicatedFieldInSpr3_negativeTest"
$ws.Range("B51:E51").HorizontalAlignment = -4108

$ws.Range("B52").Value = "rate"
$ws.Range("C52").Value = "totalValue"
$ws.Range("D52").Value = "_res_.`$MonthlyPremium(-1)"
$ws.Range("E52").Value = "_res_.`$MonthlyPremium(-2)"

$ws.Range("B53").Value = "rate"
$ws.Range("C53").Value = "totalValue"
$ws.Range("D53").Value = "_res_.`$MonthlyPremium"
$ws.Range("E53").Value = "_res_.`$MonthlyPremium"

$ws.Range("B54").Value = 0.05
$ws.Range("C54").Value = 5000
$ws.Range("C54").NumberFormat = '"$"#,##0.00'
$ws.Range("D54").Value = 250
$ws.Range("E54").Value = 240

# ---- Blank filler rows (39-42, 47-50) ----
$ws.Range("B39:E42").NumberFormat = "General"
$ws.Range("B47:E50").NumberFormat = "General"

# ---- View state ----
$ws.Range("G32").Select()
$excel.ActiveWindow.ScrollRow = 28
